$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("wasser2")

# Correct the measurement value in B8 (was 34, should be 54)
$ws.Range("B8").Value = 54

# Correct the error-propagation formula in column G: the last term under the
# square root needs to be squared as well.
$ws.Range("G2").Formula = "=SQRT((1/(B2-D2)-B2/((B2-D2)^2))^2*C2^2+(E2*B2/((B2-D2)^2))^2)"
$ws.Range("G3:G11").Formula = "=SQRT((1/(B3-D3)-B3/((B3-D3)^2))^2*C3^2+(E3*B3/((B3-D3)^2))^2)"

$wb.Application.CalculateFull()
